$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to make room for the header row
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "Supplier ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"

# Place cursor/selection on C1 to match the saved selection state
$ws.Range("C1").Select()
